$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Pre-format the numeric-looking text columns as Text so Excel stores them
# as strings (matching the source data) instead of auto-coercing to numbers.
$textCols = @('D', 'E', 'P', 'AB', 'AF', 'AJ', 'AL', 'AM')
foreach ($col in $textCols) {
    $ws.Range($col + $row).NumberFormat = "@"
}

$ws.Range('A' + $row).Value = 'Samsung-Galaxy S7 Edge-Generic'
$ws.Range('B' + $row).Value = 'DEFAULT'
$ws.Range('C' + $row).Value = 'Service_NSW], powerControl, reboot, powerSupply, [source, SERVER, offlineCharging, DEFAULT'
$ws.Range('D' + $row).Value = '2560'
$ws.Range('E' + $row).Value = '1'
$ws.Range('F' + $row).Value = 'English'
$ws.Range('G' + $row).Value = '1440x2560'
$ws.Range('H' + $row).Value = '9886783859324B4D38'
$ws.Range('I' + $row).Value = 'ac:5f:3e:2b:3d:eb'
$ws.Range('J' + $row).Value = 'Samsung'
$ws.Range('K' + $row).Value = 'Samsung-Galaxy S7 Edge'
$ws.Range('L' + $row).Value = 'OPENED'
$ws.Range('M' + $row).Value = 'SYD-L15O2-13/VIRTUAL/02'
$ws.Range('N' + $row).Value = 'srirupa.alapati@service.nsw.gov.au'
$ws.Range('O' + $row).Value = '6.0.1'
$ws.Range('P' + $row).Value = '0.01'
$ws.Range('Q' + $row).Value = 'param'
$ws.Range('R' + $row).Value = 'true'
$ws.Range('S' + $row).Value = 'srirupa.alapati@service.nsw.gov.au'
$ws.Range('T' + $row).Value = 'Galaxy S7 Edge'
$ws.Range('U' + $row).Value = 'portrait'
$ws.Range('V' + $row).Value = 'hero2ltexx-user 6.0.1 MMB29K G935FXXU1APC8 release-keys'
$ws.Range('W' + $row).Value = 'signIn'
$ws.Range('X' + $row).Value = 'srirupa.alapati@service.nsw.gov.au'
$ws.Range('Y' + $row).Value = 'Android'
$ws.Range('Z' + $row).Value = 'rotate'
$ws.Range('AA' + $row).Value = '17.02.2015'
$ws.Range('AB' + $row).Value = '0'
$ws.Range('AC' + $row).Value = 'Generic'
$ws.Range('AD' + $row).Value = '2016-06-09:01-47-45'
$ws.Range('AE' + $row).Value = 'SYD-L15O2-13/VIRTUAL/02'
$ws.Range('AF' + $row).Value = '1465436865486'
$ws.Range('AG' + $row).Value = '9886783859324B4D38'
$ws.Range('AH' + $row).Value = 'APAC-AUS-SYD'
$ws.Range('AI' + $row).Value = '2016-06-09 11:47:57'
$ws.Range('AJ' + $row).Value = '1440'
$ws.Range('AK' + $row).Value = 'pass'
$ws.Range('AL' + $row).Value = '100'
$ws.Range('AM' + $row).Value = '358809079259935'
$ws.Range('AN' + $row).Value = 'CONNECTED'
$ws.Range('AO' + $row).Value = 'mobile'
$ws.Range('AP' + $row).Value = 'C:\Users\AvoComp13\Documents\dlp-automation\test-output\screenshots\2016-06-09-11-47-57-558-AEST.png'

# New hyperlink on AP46 pointing at the screenshot for this test run
$target = 'C:/Users/AvoComp13/Documents/dlp-automation/test-output/screenshots/2016-06-09-11-47-57-558-AEST.png'
$ws.Hyperlinks.Add($ws.Range("AP" + $row), $target)
